$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Review the "OS Task Initial Release" change entry: mark it as Done
# (was "On Process") after review.
$ws.Range("E4").Value = "Done"

# Reflect the reviewer's last selection on the sheet (G4:G5).
$null = $ws.Range("G4:G5").Select()
